# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled from source for the affected rows
$ws.Range("F2").Value = 6
$ws.Range("F6").Value = -3
$ws.Range("F11").Value = -5
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = -1
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -3
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 7
$ws.Range("F31").Value = -7
$ws.Range("F33").Value = -5
$ws.Range("F44").Value = 1
$ws.Range("F45").Value = 1
$ws.Range("F47").Value = 4
$ws.Range("F57").Value = -12
$ws.Range("F60").Value = -3
